# Initial changes for 2025.
# 1) Update the cached "Date Placeholder" text (type datetimeFigureOut) from
#    10/8/23 -> 10/6/24 on the slide master and every slide layout.
# 2) Update the title text "exploreCSR 2024:" -> "exploreCSR 2025:" on slide 1.

$p = $ppt.ActivePresentation

$oldDate = "10/8/23"
$newDate = "10/6/24"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        $isDatePlaceholder = $false
        try {
            if ($shp.PlaceholderFormat.Type -eq 16) {
                $isDatePlaceholder = $true
            }
        } catch {
            $isDatePlaceholder = $false
        }
        if ($isDatePlaceholder -and $shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# Slide Master date placeholder.
$sm = $p.SlideMaster
Update-DatePlaceholder $sm.Shapes

# Every slide layout's date placeholder.
for ($li = 1; $li -le $sm.CustomLayouts.Count; $li++) {
    $cl = $sm.CustomLayouts.Item($li)
    Update-DatePlaceholder $cl.Shapes
}

# Slide 1 title text: " 2024:" -> " 2025:".
$s = $p.Slides.Item(1)
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        $tr = $shp.TextFrame.TextRange
        $full = $tr.Text
        $idx = $full.IndexOf(" 2024:")
        if ($idx -ge 0) {
            $sub = $tr.Characters($idx + 1, 6)
            $sub.Text = " 2025:"
        }
    }
}
